$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.861.73"
$ws.Range("E2").Value = "  +4.31%  "

# Row 3
$ws.Range("D3").Value = "3.769.00"
$ws.Range("E3").Value = "  +6.64%  "

# Row 4
$ws.Range("E4").Value = "  -0.28%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "426.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.24%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.87%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.22%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.737"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000309"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.31%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.64%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +15.04%  "

# Row 14
$ws.Range("D14").Value = "4.387.00"
$ws.Range("E14").Value = "  +6.96%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +12.41%  "

# Row 16
$ws.Range("E16").Value = "  +1.19%  "

# Row 17
$ws.Range("D17").Value = "3.777.93"
$ws.Range("E17").Value = "  +6.88%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.24%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +11.20%  "

# Row 20
$ws.Range("D20").Value = "66.086.04"
$ws.Range("E20").Value = "  +4.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "406.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.37%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.99%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.64%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "36.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.64%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.25%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +44.76%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.54%  "

# Row 29
$ws.Range("E29").Value = "  -0.67%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +18.05%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "696.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.11%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.128"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +17.08%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.33%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.19%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +41.45%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.150"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.70%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.97%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0473"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.06%  "

# Row 40
$ws.Range("E40").Value = "  +9.60%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.67%  "

# Row 42
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.93%  "

# Row 43
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +47.63%  "

# Row 44
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0666"
$ws.Range("E44").Value = "  +3.98%  "

# Row 45
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.12%  "

# Row 46
$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.316"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +15.96%  "

# Row 48
$ws.Range("E48").Value = "  +8.05%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.50%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.95%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.58%  "
